$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Fill in the two new Activity Log entries (rows 24 and 25) that record the
# team's report review/revision session, continuing directly after row 23.
# ---------------------------------------------------------------------------

# Column C on this sheet stores the date as text (e.g. "5-4-2020") that is
# visually re-formatted by the cell's existing date NumberFormat. Writing the
# literal string normally causes Excel to silently convert it into a real
# date serial, so we briefly force a text format, assign the literal value,
# then restore the sheet's normal date display format.
$dateFormat = $ws.Range("C23").NumberFormat

# Row 24: "Reviewed report together with team member for possible issues."
$ws.Range("B24").Value = 779
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "5-4-2020"
$ws.Range("C24").NumberFormat = $dateFormat
$ws.Range("D24").Value = 0.98402777777777783
$ws.Range("D24").Borders.LineStyle = -4142
$ws.Range("E24").Value = 0.99097222222222225
$ws.Range("G24").Value = "Reviewed report together with team member for possible issues."

# Row 25: "Revised report together with team member for clarity"
$ws.Range("B25").Value = 779
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "5-4-2020"
$ws.Range("C25").NumberFormat = $dateFormat
$ws.Range("D25").Value = 0.99097222222222225
$ws.Range("D25").Borders.LineStyle = -4142
$ws.Range("E25").Value = 0.99652777777777779
$ws.Range("G25").Value = "Revised report together with team member for clarity"

# ---------------------------------------------------------------------------
# Update the view state to match where the author left the sheet scrolled /
# zoomed to and which cell was selected.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 70
$ws.Range("D29").Select()
